# Routed ethernet and fixed ethernet symbol
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (CC0603KRX7R9BB681, 680pF) - updated unit price / subtotal
$ws.Range("G3").Value2 = 0.0107
$ws.Range("H3").Value2 = 0.0107

# Row 5 (CL32B226KAJNNNE, 22uF) - updated unit price / subtotal
$ws.Range("G5").Value2 = 0.61
$ws.Range("H5").Value2 = 1.22

# Row 6 (CL10B105MO8NNWC, 1uF) - updated unit price / subtotal
$ws.Range("G6").Value2 = 0.03374
$ws.Range("H6").Value2 = 0.3374

# Row 7 (CL10B104KA8NNNC, 100nF) - updated unit price / subtotal
$ws.Range("G7").Value2 = 0.0077
$ws.Range("H7").Value2 = 0.077

# Row 9 - ethernet symbol / part fixed: manufacturer part number changed
# from numeric 885012006044 to GRM1885C1H102JA01D, plus updated pricing
$ws.Range("C9").Value2 = "GRM1885C1H102JA01D"
$ws.Range("G9").Value2 = 0.0144
$ws.Range("H9").Value2 = 0.0144

# Row 11 (GRM21BR61E106KA73K, 10uF) - updated unit price / subtotal
$ws.Range("G11").Value2 = 0.0737
$ws.Range("H11").Value2 = 0.5159

# Row 12 (CL21B104KCFNNNE, 100nF) - updated unit price / subtotal
$ws.Range("G12").Value2 = 0.0413
$ws.Range("H12").Value2 = 0.0413

# Row 14 (CL21B473KCCWPNC, 47nF) - updated unit price / subtotal
$ws.Range("G14").Value2 = 0.0149
$ws.Range("H14").Value2 = 0.0149

# Row 18 (LQM18PN4R7MFRL, 4.7uH) - updated unit price / subtotal
$ws.Range("G18").Value2 = 0.189
$ws.Range("H18").Value2 = 0.189

# Row 20 (1SH-A-02-TS-SMT header) - ethernet header symbol fixed:
# pricing info removed (no supplier price data available anymore)
$ws.Range("G20").ClearContents()
$ws.Range("H20").ClearContents()

# Row 22 (DMG7430LFG-7) - updated unit price
$ws.Range("G22").Value2 = 0.301

# Row 30 (CRGH0603J100K, 100k) - updated unit price / subtotal
$ws.Range("G30").Value2 = 0.019
$ws.Range("H30").Value2 = 0.19

# Row 33 (73L3R10J, 100m) - updated unit price / subtotal
$ws.Range("G33").Value2 = 0.0319
$ws.Range("H33").Value2 = 0.0638

# Row 38 (RC0805FR-07100RL, 100) - updated unit price / subtotal
$ws.Range("G38").Value2 = 0.0044
$ws.Range("H38").Value2 = 0.0044
